$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended below the existing data (rows 1-18 -> now 1-20).
$ws.Range("A19").Value = "Bot yangilandi Test uchun!!!"

# Phone numbers start with "+" and must stay as literal text (not be
# coerced into a number), so mark the cell as Text before assigning it,
# then drop the explicit style again so no extra formatting lingers.
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "+998901234567"
$ws.Range("B19").Style = "Normal"

$ws.Range("A20").Value = "Sevinch"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "+998901528616"
$ws.Range("B20").Style = "Normal"
